$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 17.08155333333333
$ws.Range("N2").Value = 51.24466
$ws.Range("O2").Value = 0.3501540759902865
$ws.Range("P2").Value = 0.3501540759902865
$ws.Range("Q2").Value = 2.731346071851111
$ws.Range("R2").Value = 24.58211464666
$ws.Range("S2").Value = 0.009271281382979482
$ws.Range("T2").Value = 0.009271281382979482

# Row 3
$ws.Range("O3").Value = 0.2142771237573249
$ws.Range("P3").Value = 0.2142771237573249
$ws.Range("S3").Value = 0.005673569564116075
$ws.Range("T3").Value = 0.005673569564116076

# Row 4
$ws.Range("M4").Value = 8.398122666666666
$ws.Range("N4").Value = 25.194368
$ws.Range("O4").Value = 0.1721527793764119
$ws.Range("P4").Value = 0.1721527793764119
$ws.Range("Q4").Value = 1.342862613774222
$ws.Range("R4").Value = 12.085763523968
$ws.Range("S4").Value = 0.004558212992228537
$ws.Range("T4").Value = 0.004558212992228537

# Row 5
$ws.Range("M5").Value = 4.514486333333333
$ws.Range("N5").Value = 13.543459
$ws.Range("O5").Value = 0.09254227409953211
$ws.Range("P5").Value = 0.09254227409953213
$ws.Range("Q5").Value = 0.7218678695287778
$ws.Range("R5").Value = 6.496810825759
$ws.Range("S5").Value = 0.002450308369454415
$ws.Range("T5").Value = 0.002450308369454416

# Row 6
$ws.Range("M6").Value = 8.335727666666667
$ws.Range("N6").Value = 25.007183
$ws.Range("O6").Value = 0.1708737467764446
$ws.Range("P6").Value = 0.1708737467764446
$ws.Range("Q6").Value = 1.332885632475889
$ws.Range("R6").Value = 11.995970692283
$ws.Range("S6").Value = 0.004524347125898796
$ws.Range("T6").Value = 0.004524347125898797

# Row 7
$ws.Range("M7").Value = 17.08155333333333
$ws.Range("N7").Value = 51.24466
$ws.Range("O7").Value = 0.3501540759902865
$ws.Range("P7").Value = 0.3501540759902865
$ws.Range("Q7").Value = 100.4250484427733
$ws.Range("R7").Value = 903.82543598496
$ws.Range("S7").Value = 0.340882794607307
$ws.Range("T7").Value = 0.340882794607307

# Row 8
$ws.Range("O8").Value = 0.2142771237573249
$ws.Range("P8").Value = 0.2142771237573249
$ws.Range("S8").Value = 0.2086035541932088
$ws.Range("T8").Value = 0.2086035541932088

# Row 9
$ws.Range("M9").Value = 8.398122666666666
$ws.Range("N9").Value = 25.194368
$ws.Range("O9").Value = 0.1721527793764119
$ws.Range("P9").Value = 0.1721527793764119
$ws.Range("Q9").Value = 49.37383967197866
$ws.Range("R9").Value = 444.364557047808
$ws.Range("S9").Value = 0.1675945663841834
$ws.Range("T9").Value = 0.1675945663841834

# Row 10
$ws.Range("M10").Value = 4.514486333333333
$ws.Range("N10").Value = 13.543459
$ws.Range("O10").Value = 0.09254227409953211
$ws.Range("P10").Value = 0.09254227409953213
$ws.Range("Q10").Value = 26.54135135558933
$ws.Range("R10").Value = 238.872162200304
$ws.Range("S10").Value = 0.09009196573007769
$ws.Range("T10").Value = 0.09009196573007772

# Row 11
$ws.Range("M11").Value = 8.335727666666667
$ws.Range("N11").Value = 25.007183
$ws.Range("O11").Value = 0.1708737467764446
$ws.Range("P11").Value = 0.1708737467764446
$ws.Range("Q11").Value = 49.00700998293867
$ws.Range("R11").Value = 441.063089846448
$ws.Range("S11").Value = 0.1663493996505458
$ws.Range("T11").Value = 0.1663493996505458

